$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (well outside the used range) used to write literal text dates
# without Excel auto-converting "YYYY-MM-DD" strings into date serials.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

function Set-TextDate([string]$cellRef, [string]$text) {
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------- Row 5 ----------------
Set-TextDate "N5" "2025-10-01"
$ws.Range("Q5").Value = 42000
$ws.Range("R5").Value = -29000
$ws.Range("S5").Value = -3000
$ws.Range("T5").Value = 104000
$ws.Range("U5").Value = -23000

# ---------------- Row 29 ----------------
Set-TextDate "N29" "2025-11-05"
$ws.Range("Q29").Value = 2.21

# ---------------- Row 30 ----------------
$ws.Range("F30").Value = 0.01923414103867871
$ws.Range("G30").Value = -0.02348739337531547
Set-TextDate "N30" "2025-11-05"
$ws.Range("Q30").Value = 2.3
$ws.Range("R30").Value = 2.29
$ws.Range("S30").Value = 2.31
$ws.Range("U30").ClearContents()

# ---------------- Row 31 ----------------
$ws.Range("G31").Value = 0.03352046419083723

# ---------------- Row 47 ----------------
Set-TextDate "N47" "2025-11-04"
$ws.Range("R47").Value = 3.87
$ws.Range("U47").Value = 3.86

# ---------------- Row 48 ----------------
Set-TextDate "N48" "2025-11-04"
$ws.Range("Q48").Value = 3.58
$ws.Range("R48").Value = 3.6
$ws.Range("T48").ClearContents()
$ws.Range("U48").Value = 3.6

# ---------------- Row 49 ----------------
Set-TextDate "N49" "2025-11-04"
$ws.Range("Q49").Value = 3.69
$ws.Range("R49").Value = 3.72
$ws.Range("T49").ClearContents()
$ws.Range("U49").Value = 3.71

# ---------------- Row 50 ----------------
Set-TextDate "N50" "2025-11-04"
$ws.Range("Q50").Value = 4.1
$ws.Range("R50").Value = 4.13
$ws.Range("T50").ClearContents()

# ---------------- Row 52 ----------------
Set-TextDate "N52" "2025-11-04"
$ws.Range("Q52").Value = 5.82
$ws.Range("R52").Value = 5.84
$ws.Range("T52").ClearContents()
$ws.Range("U52").Value = 5.8

# Clean up the helper column entirely so it leaves no trace in the sheet.
$helper.EntireColumn.Delete()
